# Major reorg to study 2: add explicit `sessionEnd` (and, for Steam, also
# `sessionStart`) variables to the per-platform gaming codebooks, and a
# `sessionEnd` variable to the Android app-usage codebook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Nintendo: insert "sessionEnd" row right after "duration", before "genre"
# ---------------------------------------------------------------------
$wsNintendo = $wb.Worksheets.Item("Nintendo")
$wsNintendo.Rows.Item(9).Insert()
$wsNintendo.Cells.Item(9, 1).Value = "sessionEnd"

# ---------------------------------------------------------------------
# Xbox: same insertion as Nintendo
# ---------------------------------------------------------------------
$wsXbox = $wb.Worksheets.Item("Xbox")
$wsXbox.Rows.Item(9).Insert()
$wsXbox.Cells.Item(9, 1).Value = "sessionEnd"

# ---------------------------------------------------------------------
# Steam: insert both "sessionStart" and "sessionEnd" rows right before
# the trailing "platform" row, and widen column A to match the other
# per-platform codebook sheets.
# ---------------------------------------------------------------------
$wsSteam = $wb.Worksheets.Item("Steam")
$wsSteam.Rows.Item(11).Insert()
$wsSteam.Rows.Item(12).Insert()
$wsSteam.Cells.Item(11, 1).Value = "sessionStart"
$wsSteam.Cells.Item(12, 1).Value = "sessionEnd"
$wsSteam.Columns.Item(1).ColumnWidth = 11.83

# ---------------------------------------------------------------------
# Android: insert "sessionEnd" row right after "duration", before "category"
# ---------------------------------------------------------------------
$wsAndroid = $wb.Worksheets.Item("Android")
$wsAndroid.Rows.Item(10).Insert()
$wsAndroid.Cells.Item(10, 1).Value = "sessionEnd"
